$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '257.24' },
    @{ Cell = 'E2'; Value = '0.86%' },
    @{ Cell = 'D3'; Value = '27.08' },
    @{ Cell = 'E3'; Value = '-3.56%' },
    @{ Cell = 'D4'; Value = '4.763' },
    @{ Cell = 'E4'; Value = '-10.04%' },
    @{ Cell = 'D5'; Value = '0.05948' },
    @{ Cell = 'E5'; Value = '1.65%' },
    @{ Cell = 'D6'; Value = '6.654' },
    @{ Cell = 'E6'; Value = '-0.73%' },
    @{ Cell = 'D7'; Value = '0.8725' },
    @{ Cell = 'E7'; Value = '0.28%' },
    @{ Cell = 'D8'; Value = '0.9533' },
    @{ Cell = 'E8'; Value = '2.71%' },
    @{ Cell = 'D9'; Value = '0.1406' },
    @{ Cell = 'E9'; Value = '-0.41%' },
    @{ Cell = 'D10'; Value = '0.03902' },
    @{ Cell = 'E10'; Value = '12.78%' },
    @{ Cell = 'D11'; Value = '0.07170' },
    @{ Cell = 'E11'; Value = '0.95%' },
    @{ Cell = 'D12'; Value = '0.03197' },
    @{ Cell = 'E12'; Value = '0.65%' },
    @{ Cell = 'D13'; Value = '0.09270' },
    @{ Cell = 'E13'; Value = '0.49%' },
    @{ Cell = 'D14'; Value = '0.001540' },
    @{ Cell = 'E14'; Value = '-0.72%' },
    @{ Cell = 'D15'; Value = '0.0006060' },
    @{ Cell = 'E15'; Value = '0.34%' },
    @{ Cell = 'D16'; Value = '0.006046' },
    @{ Cell = 'E16'; Value = '4.20%' },
    @{ Cell = 'D17'; Value = '3.482' },
    @{ Cell = 'E17'; Value = '-0.44%' },
    @{ Cell = 'D18'; Value = '3.196' },
    @{ Cell = 'E18'; Value = '-1.10%' },
    @{ Cell = 'E19'; Value = '1.68%' },
    @{ Cell = 'E20'; Value = '-1.44%' },
    @{ Cell = 'D21'; Value = '0.1306' },
    @{ Cell = 'E21'; Value = '-0.73%' },
    @{ Cell = 'D22'; Value = '3.812' },
    @{ Cell = 'E22'; Value = '8.01%' },
    @{ Cell = 'E23'; Value = '1.29%' },
    @{ Cell = 'D25'; Value = '0.001226' },
    @{ Cell = 'E25'; Value = '-0.40%' },
    @{ Cell = 'D26'; Value = '0.004501' },
    @{ Cell = 'E26'; Value = '-10.01%' },
    @{ Cell = 'E27'; Value = '0.01%' },
    @{ Cell = 'D28'; Value = '0.0001493' },
    @{ Cell = 'E28'; Value = '86.65%' },
    @{ Cell = 'D40'; Value = '0.03833' },
    @{ Cell = 'E40'; Value = '0.37%' },
    @{ Cell = 'D41'; Value = '0.006180' },
    @{ Cell = 'E41'; Value = '62.04%' },
    @{ Cell = 'E42'; Value = '-0.07%' },
    @{ Cell = 'E43'; Value = '-4.49%' },
    @{ Cell = 'D44'; Value = '0.01057' },
    @{ Cell = 'E44'; Value = '4.85%' },
    @{ Cell = 'D45'; Value = '0.00005500' },
    @{ Cell = 'E45'; Value = '5.26%' },
    @{ Cell = 'E46'; Value = '0.02%' },
    @{ Cell = 'D47'; Value = '0.08852' },
    @{ Cell = 'D48'; Value = '0.002386' },
    @{ Cell = 'E48'; Value = '10.78%' },
    @{ Cell = 'E49'; Value = '0.02%' },
    @{ Cell = 'E50'; Value = '0.02%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
